$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text storage so
# values such as "515.67" are not auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 34/35: NEARProtocol and EnergySwap swap places with updated stats
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "4.43"
$ws.Range("E34").Value = "  -2.64%  "

$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").Value = "26.88"
$ws.Range("E35").Value = "  +7.84%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
$updates = @{
    2 = @("57.307.37", "  +2.61%  ")
    3 = @("3.058.66", "  +3.83%  ")
    4 = @($null, "  +0.88%  ")
    5 = @("515.67", "  +3.17%  ")
    6 = @("139.63", "  +2.32%  ")
    7 = @($null, "  +0.16%  ")
    8 = @("0.432", "  +0.94%  ")
    9 = @("7.27", "  +2.93%  ")
    10 = @($null, "  +1.83%  ")
    11 = @("0.372", "  +1.93%  ")
    12 = @("3.591.11", "  +2.82%  ")
    13 = @($null, "  +3.45%  ")
    14 = @("25.46", "  -2.26%  ")
    15 = @("0.0000162", "  +0.76%  ")
    16 = @("57.395.03", "  +3.30%  ")
    17 = @("3.048.18", "  +2.61%  ")
    18 = @($null, "  +1.94%  ")
    19 = @("12.90", "  -0.04%  ")
    20 = @("8.07", "  +3.04%  ")
    21 = @("332.67", "  +2.39%  ")
    22 = @($null, "  +0.70%  ")
    23 = @("0.500", "  +1.73%  ")
    24 = @("65.69", "  +2.96%  ")
    25 = @("0.169", "  +5.54%  ")
    26 = @("1.00", "  -0.86%  ")
    27 = @("0.0₃0904", "  +2.78%  ")
    28 = @("6.31", "  -2.61%  ")
    29 = @("7.17", "  +2.27%  ")
    30 = @("1.81", "  +2.30%  ")
    31 = @("20.74", "  +2.97%  ")
    32 = @($null, "  +0.59%  ")
    33 = @("154.14", "  +2.22%  ")
    36 = @("5.87", "  +1.71%  ")
    37 = @($null, "  +1.35%  ")
    38 = @("0.0671", "  +2.15%  ")
    39 = @("3.099.65", "  +3.80%  ")
    40 = @("3.90", "  +3.96%  ")
    41 = @("36.84", "  +0.71%  ")
    42 = @($null, "  +1.28%  ")
    43 = @("0.653", "  +1.57%  ")
    44 = @("2.260.13", "  +4.72%  ")
    45 = @("0.0255", "  +8.29%  ")
    46 = @("1.37", "  +1.39%  ")
    47 = @("20.31", "  +3.49%  ")
    48 = @("5.85", "  -0.33%  ")
    49 = @("0.923", "  -1.06%  ")
    50 = @("261.57", "  +15.42%  ")
    51 = @("0.0871", "  +2.61%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals[0]) {
        $ws.Cells.Item($row, 4).Value = $vals[0]
    }
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

Write-Output "Applied cryptos update"
